$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper note: D-column prices are text-formatted (e.g. "23.215.85", "0.9982")
# in this sheet, so new values are entered with a leading quote to force text
# entry (like typing '0.9993 in Excel) and the style is reset to Normal so
# the quote-prefix flag does not leave a stray format change on the cell.

# Rows 2-10: refreshed Price / Volume(1h) figures
$ws.Range("D2").Value = "'23.194.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "'1.602.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'0.9988"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'302.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").Value = "'0.3778"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'52.08"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.20%  "
$ws.Range("D9").Value = "'0.3612"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").Value = "'1.265"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.22%  "

# Rows 11-12: Dogecoin and BinanceUSD swapped order, with refreshed figures
$ws.Range("B11").Value = "BinanceUSD"
$ws.Range("C11").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D11").Value = "'0.9989"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "'0.08122"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.43%  "

# Rows 13-51: refreshed Price / Volume(1h) figures
$ws.Range("D13").Value = "'22.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("D14").Value = "'6.603"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "'7.397"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "'0.00001251"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").Value = "'1.603.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "'93.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").Value = "'0.06860"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "'18.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").Value = "'0.9984"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Value = "'23.186.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "'2.395"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.52%  "
$ws.Range("D26").Value = "'3.002"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.95%  "
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").Value = "'149.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").Value = "'5.228"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("D30").Value = "'134.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").Value = "'2.418"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("D32").Value = "'6.819"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("D33").Value = "'1.781.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").Value = "'0.9836"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.87%  "
$ws.Range("D35").Value = "'0.07598"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("D36").Value = "'10.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.34%  "
$ws.Range("D37").Value = "'0.02723"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").Value = "'6.155"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("D39").Value = "'0.2506"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("D40").Value = "'0.08791"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("D41").Value = "'0.7125"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "'1.362"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.23%  "
$ws.Range("D43").Value = "'12.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("D44").Value = "'15.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").Value = "'0.6572"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("D46").Value = "'2.312"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("D48").Value = "'132.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").Value = "'0.07963"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").Value = "'1.206"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").Value = "'1.223"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.10%  "
